# Add a table of widget data below the existing invoice header/data rows.
# Row 3 is intentionally left blank (matches the target diff's sparse rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 4)
$ws.Range("A4").Value = "Item"
$ws.Range("B4").Value = "Quantity"
$ws.Range("C4").Value = "Price"
$ws.Range("D4").Value = "Amount"

# Data rows (5-10)
$ws.Range("A5").Value = "Widget A"
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = 25.5
$ws.Range("D5").Value = 255

$ws.Range("A6").Value = "Widget B"
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 40
$ws.Range("D6").Value = 200

$ws.Range("A7").Value = "Widget C"
$ws.Range("B7").Value = 15
$ws.Range("C7").Value = 30
$ws.Range("D7").Value = 450

$ws.Range("A8").Value = "Widget D"
$ws.Range("B8").Value = 8
$ws.Range("C8").Value = 50
$ws.Range("D8").Value = 400

$ws.Range("A9").Value = "Widget E"
$ws.Range("B9").Value = 12
$ws.Range("C9").Value = 35
$ws.Range("D9").Value = 420

$ws.Range("A10").Value = "Widget F"
$ws.Range("B10").Value = 20
$ws.Range("C10").Value = 28
$ws.Range("D10").Value = 560
